$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-4 (replacing old rows 2,3,4 and dropping old rows 5,6,7)
# Row 2: ECs -> FAPs
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.1716463333333333
$ws.Range("H2").Value = 0.514939
$ws.Range("I2").Value = 0.07772289907851986
$ws.Range("J2").Value = 0.07772289907851984
$ws.Range("M2").Value = 0.001809666666666667
$ws.Range("N2").Value = 0.005429
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0003106226478888889
$ws.Range("R2").Value = 0.002795603831
$ws.Range("S2").Value = 0.07772289907851986
$ws.Range("T2").Value = 0.07772289907851984

# Row 3: FAPs -> ECs
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "ECs"
$ws.Range("G3").Value = 1.440925666666667
$ws.Range("H3").Value = 4.322777
$ws.Range("I3").Value = 0.652463224789629
$ws.Range("J3").Value = 0.652463224789629
$ws.Range("M3").Value = 5.973506333333333
$ws.Range("N3").Value = 17.920519
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.002607595148111111
$ws.Range("R3").Value = 0.023468356333
$ws.Range("S3").Value = 0.652463224789629
$ws.Range("T3").Value = 0.652463224789629

# Row 4: MuSCs -> FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 0.5958676666666666
$ws.Range("H4").Value = 1.787603
$ws.Range("I4").Value = 0.2698138761318511
$ws.Range("J4").Value = 0.2698138761318511
$ws.Range("M4").Value = 0.001809666666666667
$ws.Range("N4").Value = 0.005429
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.001078321854111111
$ws.Range("R4").Value = 0.009704896686999998
$ws.Range("S4").Value = 0.2698138761318511
$ws.Range("T4").Value = 0.2698138761318511

# Remove old rows 5-7 (data no longer present)
$ws.Range("A5:T7").EntireRow.Delete()
